# Generate Report for Handback
# Refresh the timestamp cells on the handback-status report to reflect the
# latest run (Overview "Latest HO Xliff Generate Date" plus the per-locale
# "Correspond Handoff/Handback Datetime" columns for the first data row).

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview sheet: "Latest HO Xliff Generate Date" for row 2 (en mirrors the
# de-de handoff datetime, since de-de is generated last for this file).
$wsOverview.Range("G2").Value = "2016-08-18 13:06:20"

# zh-cn sheet: Correspond Handoff Datetime / Correspond Handback DateTime
$wsZhCn.Range("H2").Value = "2016-08-18 13:06:14"
$wsZhCn.Range("K2").Value = "2016-08-18 13:06:31"

# de-de sheet: Correspond Handoff Datetime / Correspond Handback DateTime
$wsDeDe.Range("H2").Value = "2016-08-18 13:06:20"
$wsDeDe.Range("K2").Value = "2016-08-18 13:06:39"
